$wb = $excel.ActiveWorkbook
$wsSettings = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")

# --- Settings sheet: delete blank row 6 (shift rows 7+ up by one) ---
$wsSettings.Rows.Item(6).Delete()

# --- Constants sheet: update MaxRetryNumber value and add Statuses row ---
$wsConstants.Range("B2").Value = 2
$wsConstants.Range("A20").Value = "Statuses"
$wsConstants.Range("B20").Value = "Retrieved, Added, Installed, Failed"
